$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its text formatting so numeric-looking
# strings (e.g. "241.04", "1.980.26") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.942.98"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "1.821.41"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("D5").Value = "241.04"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  -5.38%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "44.29"
$ws.Range("E8").Value = "  +5.14%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.07297"
$ws.Range("E9").Value = "  -1.93%  "
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "0.2911"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("D11").Value = "22.73"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").Value = "0.07649"
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("D13").Value = "1.829.54"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").Value = "0.6605"
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("D16").Value = "82.05"
$ws.Range("E16").Value = "  -4.54%  "
$ws.Range("D17").Value = "6.037"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("D18").Value = "0.000008569"
$ws.Range("E18").Value = "  +2.89%  "
$ws.Range("D19").Value = "28.964.19"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("D20").Value = "2.082.48"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("D22").Value = "222.88"
$ws.Range("E22").Value = "  -2.67%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "7.050"
$ws.Range("E24").Value = "  -1.94%  "
$ws.Range("D25").Value = "1.001"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "157.81"
$ws.Range("E26").Value = "  -2.07%  "
$ws.Range("D27").Value = "8.411"
$ws.Range("E27").Value = "  -3.49%  "
$ws.Range("D28").Value = "0.1366"
$ws.Range("E28").Value = "  -2.90%  "
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").Value = "1.498"
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("D31").Value = "4.065"
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("D33").Value = "3.993"
$ws.Range("E33").Value = "  -2.04%  "
$ws.Range("D34").Value = "0.05275"
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("D35").Value = "0.7364"
$ws.Range("E35").Value = "  -2.88%  "
$ws.Range("D36").Value = "1.821"
$ws.Range("E36").Value = "  -2.69%  "
$ws.Range("D37").Value = "1.148"
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("D38").Value = "2.646"
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("D39").Value = "1.288.75"
$ws.Range("E39").Value = "  -2.61%  "
$ws.Range("D40").Value = "2.737"
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("D41").Value = "0.01773"
$ws.Range("E41").Value = "  -1.70%  "
$ws.Range("D42").Value = "6.327"
$ws.Range("E42").Value = "  +5.59%  "
$ws.Range("D43").Value = "0.8924"
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "102.15"
$ws.Range("E45").Value = "  -1.17%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.00000000125"
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.980.26"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("D48").Value = "0.5137"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("D49").Value = "63.90"
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("E50").Value = "  -3.14%  "
$ws.Range("D51").Value = "0.07334"
$ws.Range("E51").Value = "  -11.66%  "
